# Insert a new data row at row 103 (pushes the existing rows 103-209 down to
# 104-210, matching the canonical diff). Excel's Insert() defaults to
# shifting cells down and carries formatting from the row above, which is
# exactly what the target XML shows (only column D keeps style index 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Insert()

# Populate the newly inserted row with the new "Albahaca" price record.
# Columns A, B, C, E, F, G, H, I, R are unchanged vs. the row that used to
# occupy position 103 (now at 104); only D, J, K, L, M, N, O, P, Q differ.
$ws.Cells.Item(103, 1).Value2  = 9
$ws.Cells.Item(103, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(103, 3).Value2  = "Metropolitana"
$ws.Cells.Item(103, 4).Value2  = 44494
$ws.Cells.Item(103, 5).Value2  = 13
$ws.Cells.Item(103, 6).Value2  = 100112052
$ws.Cells.Item(103, 7).Value2  = "Albahaca"
$ws.Cells.Item(103, 8).Value2  = "Sin especificar"
$ws.Cells.Item(103, 9).Value2  = "Primera"
$ws.Cells.Item(103, 10).Value2 = 160
$ws.Cells.Item(103, 11).Value2 = 4000
$ws.Cells.Item(103, 12).Value2 = 4000
$ws.Cells.Item(103, 13).Value2 = 4000
$ws.Cells.Item(103, 14).Value2 = "$/paquete"
$ws.Cells.Item(103, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(103, 16).Value2 = 4000
$ws.Cells.Item(103, 17).Value2 = 1
$ws.Cells.Item(103, 18).Value2 = "Hortaliza"
